$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Kyrgyz title in A1 - wording tweak
# "Коопсуз" -> "Коопсуздук" (x2) and "жоктугунана" -> "жоктугунан"
$ws.Range("A1").Value = "3.9.2 Коопсуздук суунун, коопсуздук санитариянын жана гигиенанын жоктугунан болгон өлүм"

# Add the new 2022 data column (S), one year to the right of the existing 2021 column (R)
$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 1.2
$ws.Range("S6").Value = 2.7
$ws.Range("S7").Value = 0.9
$ws.Range("S8").Value = 0.4
$ws.Range("S9").Value = 0.7
$ws.Range("S10").Value = 0.9
$ws.Range("S11").Value = 1.1000000000000001
$ws.Range("S12").Value = 2.7
$ws.Range("S13").Value = 0.4
$ws.Range("S14").Value = 0.6

# Copy the number-format styling from the existing 2021 column (R) onto the new 2022 column (S)
$ws.Range("S4:S14").NumberFormat = $ws.Range("R4:R14").NumberFormat

$ws.Range("A1").Select()
